$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 106, pushing the existing rows 106-109 down to 107-110.
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with the new weekly price entry.
# (Same market/category metadata as the surrounding rows, new date + price.)
$ws.Cells.Item(106, 1).Value = 5
$ws.Cells.Item(106, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(106, 3).Value = "Maule"
$ws.Cells.Item(106, 4).Value = 45147
$ws.Cells.Item(106, 5).Value = 7
$ws.Cells.Item(106, 6).Value = 100112040
$ws.Cells.Item(106, 7).Value = "Cilantro"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 200
$ws.Cells.Item(106, 11).Value = 8000
$ws.Cells.Item(106, 12).Value = 8000
$ws.Cells.Item(106, 13).Value = 8000
$ws.Cells.Item(106, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(106, 15).Value = "Región Metropolitana"
$ws.Cells.Item(106, 16).Value = 222
$ws.Cells.Item(106, 17).Value = 36
$ws.Cells.Item(106, 18).Value = "Hortaliza"

# Row 108 (old 107) origin flips from "Región del Maule" to "Región Metropolitana".
$ws.Cells.Item(108, 15).Value = "Región Metropolitana"

# Row 109 (old 108) origin flips from "Región Metropolitana" to "Región del Maule".
$ws.Cells.Item(109, 15).Value = "Región del Maule"
